# "fix gantt from files"
# The task table (Table1, A1:H6) is re-sorted by EstimatedEffortHours
# (column E) descending, and the sheet selection is changed from the
# single cell B10 to the whole table range A1:H6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Excel enumerations used below (values match the real Excel object model):
#   xlSortOnValues = 0
#   xlDescending   = 2
#   xlYes          = 1  (table/range has headers)
$xlSortOnValues = 0
$xlDescending = 2
$xlYes = 1

# Sort the table by EstimatedEffortHours (column E) descending.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("E1:E6"), $xlSortOnValues, $xlDescending) | Out-Null
$lo.Sort.Header = $xlYes
$lo.Sort.Apply() | Out-Null

# Update the selection to span the whole table instead of cell B10.
$ws.Range("A1:H6").Select() | Out-Null
